# Adds two new columns, I (I0) and J (IF), to the sheet: header labels in
# row 1 (styled like the existing headers) plus per-row numeric data for
# rows 2-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting used by the other header cells (bold/centered/bordered)
# by copying the existing H1 header format onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2-37: [row, I value, J value]
$data = @(
    @(2, 7, 8),
    @(3, 8, 8),
    @(4, 7, 8),
    @(5, 6, 6),
    @(6, 6, 6),
    @(7, 3, 4),
    @(8, 7, 7),
    @(9, 6, 7),
    @(10, 7, 7),
    @(11, 6, 6),
    @(12, 7, 8),
    @(13, 6, 7),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 10, 10),
    @(17, 5, 6),
    @(18, 7, 7),
    @(19, 6, 7),
    @(20, 7, 7),
    @(21, 6, 6),
    @(22, 7, 7),
    @(23, 2, 3),
    @(24, 7, 7),
    @(25, 6, 6),
    @(26, 5, 6),
    @(27, 5, 5),
    @(28, 6, 6),
    @(29, 8, 8),
    @(30, 6, 6),
    @(31, 7, 7),
    @(32, 7, 7),
    @(33, 6, 6),
    @(34, 3, 3),
    @(35, 9, 9),
    @(36, 5, 5),
    @(37, 4, 4)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
